$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete row 264 (NNG, Nanning, China) - all subsequent rows shift up by one
$ws.Rows.Item(264).Delete()
